$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 13.03.2024"

$ws.Range("B6").Value = "14.03."
$ws.Range("C6").Value = "15.03."
$ws.Range("D6").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E6").Value = "54,03-"

$ws.Range("B7").Value = "18.03."
$ws.Range("C7").Value = "19.03."
$ws.Range("D7").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E7").Value = "56,87-"

$ws.Range("B8").Value = "19.03."
$ws.Range("C8").Value = "20.03."
$ws.Range("D8").Value = "PAYPAL SJTQOM"
$ws.Range("E8").Value = "5,70-"

$ws.Range("B9").Value = "21.03."
$ws.Range("C9").Value = "22.03."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 69197396"
$ws.Range("E9").Value = "38,83-"

$ws.Range("B10").Value = "23.03."
$ws.Range("C10").Value = "24.03."
$ws.Range("D10").Value = "MCDONALDS Sebnitz"
$ws.Range("E10").Value = "44,00-"

$ws.Range("D12").Value = "KONTOSTAND AM 28.03.2024"
$ws.Range("E12").Value = "199,43-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 06.04.2024"
